$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)  # 总计
$ws2 = $wb.Worksheets.Item(2)  # current 2022-Q2 (will remain at position 3 after insert)

function SetText($cell, $text) {
  $cell.Value = "'" + $text
  $cell.ClearFormats()
}

# --- 1) Update the "总计" (summary) sheet: insert a new top data row for 2022-Q3,
#        shifting 2022-Q2 / 2022-Q1 / 2021-Q4 rows down by one and bumping their
#        index column (A) by 1. Avoid Rows().Insert() so no stray unused styles
#        get created; copy values forward manually instead. ---

$b2 = $ws1.Cells.Item(2,2).Value2
$c2 = $ws1.Cells.Item(2,3).Value2
$d2 = $ws1.Cells.Item(2,4).Value2

$b3 = $ws1.Cells.Item(3,2).Value2
$c3 = $ws1.Cells.Item(3,3).Value2
$d3 = $ws1.Cells.Item(3,4).Value2

$b4 = $ws1.Cells.Item(4,2).Value2
$c4 = $ws1.Cells.Item(4,3).Value2
$d4 = $ws1.Cells.Item(4,4).Value2

# old row2 (2022-Q2) -> row3
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = $b2
$ws1.Cells.Item(3,3).Value = $c2
$ws1.Cells.Item(3,4).Value = $d2

# old row3 (2022-Q1) -> row4
$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = $b3
$ws1.Cells.Item(4,3).Value = $c3
$ws1.Cells.Item(4,4).Value = $d3

# old row4 (2021-Q4) -> row5 (brand-new row; copy column-A style before overwriting row4)
$ws1.Cells.Item(4,1).Copy($ws1.Cells.Item(5,1))
$ws1.Cells.Item(5,1).Value = 3
$ws1.Cells.Item(5,2).Value = $b4
$ws1.Cells.Item(5,3).Value = $c4
$ws1.Cells.Item(5,4).Value = $d4

# new row2: 2022-Q3 summary data
$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q3"
$ws1.Cells.Item(2,3).Value = 12
$ws1.Cells.Item(2,4).Value = 0.8100000000000001

# --- 2) Insert a brand-new "2022-Q3" sheet right after "总计" (i.e. before the
#        sheet currently named "2022-Q2"), with the quarter's fund holdings. ---

$newWs = $wb.Worksheets.Add($ws2)
$newWs.Name = "2022-Q3"

# Header row, reusing the bold/bordered header style already used elsewhere.
$ws1.Cells.Item(1,2).Copy($newWs.Range("B1:H1"))
$newWs.Cells.Item(1,2).Value = "基金代码"
$newWs.Cells.Item(1,3).Value = "基金名称"
$newWs.Cells.Item(1,4).Value = "基金规模"
$newWs.Cells.Item(1,5).Value = "股票总仓位"
$newWs.Cells.Item(1,6).Value = "仓位占比"
$newWs.Cells.Item(1,7).Value = "持有市值(亿元)"
$newWs.Cells.Item(1,8).Value = "仓位排名"

# Column A (row index 0..11), reusing the same style.
$ws1.Cells.Item(2,1).Copy($newWs.Range("A2:A13"))
for ($r=2; $r -le 13; $r++) {
  $newWs.Cells.Item($r,1).Value = $r - 2
}

SetText $newWs.Cells.Item(2,2) "004616"
SetText $newWs.Cells.Item(2,3) "中欧电子信息产业沪港深股票A"
SetText $newWs.Cells.Item(2,4) "5.01"
SetText $newWs.Cells.Item(2,5) "92.97"
SetText $newWs.Cells.Item(2,6) "5.34"
SetText $newWs.Cells.Item(2,7) "0.2675"
$newWs.Cells.Item(2,8).Value = 5
SetText $newWs.Cells.Item(3,2) "005763"
SetText $newWs.Cells.Item(3,3) "中欧电子信息产业沪港深股票C"
SetText $newWs.Cells.Item(3,4) "3.88"
SetText $newWs.Cells.Item(3,5) "92.97"
SetText $newWs.Cells.Item(3,6) "5.34"
SetText $newWs.Cells.Item(3,7) "0.2072"
$newWs.Cells.Item(3,8).Value = 5
SetText $newWs.Cells.Item(4,2) "001411"
SetText $newWs.Cells.Item(4,3) "诺安创新驱动灵活配置混合A"
SetText $newWs.Cells.Item(4,4) "3.98"
SetText $newWs.Cells.Item(4,5) "80.56"
SetText $newWs.Cells.Item(4,6) "2.63"
SetText $newWs.Cells.Item(4,7) "0.1047"
$newWs.Cells.Item(4,8).Value = 9
SetText $newWs.Cells.Item(5,2) "010824"
SetText $newWs.Cells.Item(5,3) "天弘创新成长混合A"
SetText $newWs.Cells.Item(5,4) "2.30"
SetText $newWs.Cells.Item(5,5) "79.97"
SetText $newWs.Cells.Item(5,6) "3.58"
SetText $newWs.Cells.Item(5,7) "0.0823"
$newWs.Cells.Item(5,8).Value = 5
SetText $newWs.Cells.Item(6,2) "005310"
SetText $newWs.Cells.Item(6,3) "广发电子信息传媒股票A"
SetText $newWs.Cells.Item(6,4) "1.55"
SetText $newWs.Cells.Item(6,5) "89.36"
SetText $newWs.Cells.Item(6,6) "3.66"
SetText $newWs.Cells.Item(6,7) "0.0567"
$newWs.Cells.Item(6,8).Value = 4
SetText $newWs.Cells.Item(7,2) "002051"
SetText $newWs.Cells.Item(7,3) "诺安创新驱动灵活配置混合C"
SetText $newWs.Cells.Item(7,4) "1.80"
SetText $newWs.Cells.Item(7,5) "80.56"
SetText $newWs.Cells.Item(7,6) "2.63"
SetText $newWs.Cells.Item(7,7) "0.0473"
$newWs.Cells.Item(7,8).Value = 9
SetText $newWs.Cells.Item(8,2) "010825"
SetText $newWs.Cells.Item(8,3) "天弘创新成长混合C"
SetText $newWs.Cells.Item(8,4) "0.92"
SetText $newWs.Cells.Item(8,5) "79.97"
SetText $newWs.Cells.Item(8,6) "3.58"
SetText $newWs.Cells.Item(8,7) "0.0329"
$newWs.Cells.Item(8,8).Value = 5
SetText $newWs.Cells.Item(9,2) "011214"
SetText $newWs.Cells.Item(9,3) "招商惠润一年定期开放混合（MOM）A"
SetText $newWs.Cells.Item(9,4) "0.48"
SetText $newWs.Cells.Item(9,5) "57.54"
SetText $newWs.Cells.Item(9,6) "1.69"
SetText $newWs.Cells.Item(9,7) "0.0081"
$newWs.Cells.Item(9,8).Value = 10
SetText $newWs.Cells.Item(10,2) "010236"
SetText $newWs.Cells.Item(10,3) "广发电子信息传媒股票C"
SetText $newWs.Cells.Item(10,4) "0.13"
SetText $newWs.Cells.Item(10,5) "89.36"
SetText $newWs.Cells.Item(10,6) "3.66"
SetText $newWs.Cells.Item(10,7) "0.0048"
$newWs.Cells.Item(10,8).Value = 4
SetText $newWs.Cells.Item(11,2) "001914"
SetText $newWs.Cells.Item(11,3) "中信建投聚利混合A"
SetText $newWs.Cells.Item(11,4) "0.10"
SetText $newWs.Cells.Item(11,5) "39.73"
SetText $newWs.Cells.Item(11,6) "2.11"
SetText $newWs.Cells.Item(11,7) "0.0021"
$newWs.Cells.Item(11,8).Value = 5
SetText $newWs.Cells.Item(12,2) "011215"
SetText $newWs.Cells.Item(12,3) "招商惠润一年定期开放混合（MOM）C"
SetText $newWs.Cells.Item(12,4) "0.06"
SetText $newWs.Cells.Item(12,5) "57.54"
SetText $newWs.Cells.Item(12,6) "1.69"
SetText $newWs.Cells.Item(12,7) "0.0010"
$newWs.Cells.Item(12,8).Value = 10
SetText $newWs.Cells.Item(13,2) "006845"
SetText $newWs.Cells.Item(13,3) "中信建投聚利混合C"
SetText $newWs.Cells.Item(13,4) "0.01"
SetText $newWs.Cells.Item(13,5) "39.73"
SetText $newWs.Cells.Item(13,6) "2.11"
SetText $newWs.Cells.Item(13,7) "0.0002"
$newWs.Cells.Item(13,8).Value = 5

